# New weekly price report: insert a new record as the first row of the
# "Mango" price series (row 73), pushing all subsequent rows (old 73-189)
# down by one (new 74-190).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 73; Excel shifts existing rows 73-189
# down to 74-190 and extends the used range to A1:T190 automatically.
$ws.Rows("73").Insert()

# Populate the newly inserted row 73 with the new weekly record.
$ws.Range("A73").Value = 4
$ws.Range("B73").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C73").Value = "Los Lagos"
$ws.Range("D73").Value = 44665
$ws.Range("E73").Value = 10
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100108
$ws.Range("H73").Value = "Tropicales y subtropicales"
$ws.Range("I73").Value = 100108002
$ws.Range("J73").Value = "Mango"
$ws.Range("K73").Value = "Sin especificar"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 200
$ws.Range("N73").Value = 7500
$ws.Range("O73").Value = 8000
$ws.Range("P73").Value = 7750
$ws.Range("Q73").Value = "$/bandeja 4 kilos"
$ws.Range("R73").Value = "Perú"
$ws.Range("S73").Value = 1938
$ws.Range("T73").Value = 4
